$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update the hotel-name / city values for the existing rows.
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("A3").Value = "London"
$ws.Range("B3").Value = "Grand Plaza Apartments"

# Move the selection to D11, matching the saved view state.
$ws.Range("D11").Select()
